# Update the "Date" footer placeholder (datetimeFigureOut field) from
# 15/05/2020 to 16/05/2020 across the slide master and every slide layout.
#
# ppPlaceholderDate = 16, msoPlaceholder = 14.

$p = $ppt.ActivePresentation
$newDate = "16/05/2020"

# --- Slide Master -----------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Type -eq 14) {
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every Slide Layout -------------------------------------------------
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
